# ---------------------------------------------------------------------------
# Edit summary (from the authoritative diff):
#   1. The table on slide 5 (the "B1 - types of financial documents" table)
#      switches from the deck's custom table style
#      {5633DA39-EF30-4012-96F7-5E03EDB2311B} to the built-in style
#      {549E92F9-1CB7-40D3-AFD6-944B548E9DEA}.
#   2. The presentation's theme colours are swapped: the slide master's
#      theme (the one driving the deck's visible design) moves from the
#      "Integral" / Red Violet palette to the default "Office Theme" /
#      Office palette.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style swap on slide 5.
# ---------------------------------------------------------------------------
$s   = $p.Slides.Item(5)
$sh  = $s.Shapes.Item(2)          # the graphicFrame holding the table
$tbl = $sh.Table
$tbl.ApplyStyle("{549E92F9-1CB7-40D3-AFD6-944B548E9DEA}")

# ---------------------------------------------------------------------------
# 2) Swap the active design's colour scheme from "Integral" (Red Violet)
#    back to the stock "Office Theme" (Office) palette.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$tcs    = $master.Theme.ThemeColorScheme

# Index order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
